# Regenerate the handback-status report's timestamp columns for the
# "47634a4a-64d0-49e2-baaf-004ec3d0286f.md" file row, as produced by a
# fresh "Generate Report" run for the handback.
$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value2 = "2017-02-21 05:29:36"

# --- zh-cn sheet: Correspond Handoff/Handback datetimes (row 2) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value2 = "2017-02-21 05:29:20"
$wsZhCn.Range("L2").Value2 = "2017-02-21 05:30:16"

# --- de-de sheet: Correspond Handoff/Handback datetimes (row 2) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value2 = "2017-02-21 05:29:36"
$wsDeDe.Range("L2").Value2 = "2017-02-21 05:30:39"
